# Adds a new "2022-Q3" sheet (with fund holding data) right after "总计",
# ahead of the existing "2022-Q2" sheet, and inserts a matching summary
# row at the top of the "总计" sheet's data table (shifting every other
# quarter's row down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet by duplicating "2022-Q2" (so it
#    inherits identical formatting: bold/bordered header row + index
#    column) and placing the copy immediately before it.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Extend the formatted index column (A) down to rows 4 and 5 so every
# data row matches the existing style (bold, centered, bordered).
$q3.Range("A2").Copy()
$q3.Range("A4:A5").PasteSpecial(-4122)

# Columns B-G hold text values (fund codes/names/percentages stored as
# strings, matching the source data) - force text type so leading zeros
# and the literal decimal formatting survive, then drop the number
# format override again so the cell settles back to the plain (no
# explicit style) look used throughout the rest of the table.
function Set-TextCell($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
$q3.Range("A2").Value = 0
Set-TextCell $q3.Range("B2") "002210"
Set-TextCell $q3.Range("C2") "创金合信量化多因子股票A"
Set-TextCell $q3.Range("D2") "2.39"
Set-TextCell $q3.Range("E2") "91.71"
Set-TextCell $q3.Range("F2") "1.22"
Set-TextCell $q3.Range("G2") "0.0292"
$q3.Range("H2").Value = 10

# Row 3
$q3.Range("A3").Value = 1
Set-TextCell $q3.Range("B3") "004194"
Set-TextCell $q3.Range("C3") "招商中证1000指数增强A"
Set-TextCell $q3.Range("D3") "1.56"
Set-TextCell $q3.Range("E3") "92.06"
Set-TextCell $q3.Range("F3") "1.04"
Set-TextCell $q3.Range("G3") "0.0162"
$q3.Range("H3").Value = 9

# Row 4
$q3.Range("A4").Value = 2
Set-TextCell $q3.Range("B4") "004195"
Set-TextCell $q3.Range("C4") "招商中证1000指数增强C"
Set-TextCell $q3.Range("D4") "1.09"
Set-TextCell $q3.Range("E4") "92.06"
Set-TextCell $q3.Range("F4") "1.04"
Set-TextCell $q3.Range("G4") "0.0113"
$q3.Range("H4").Value = 9

# Row 5
$q3.Range("A5").Value = 3
Set-TextCell $q3.Range("B5") "003865"
Set-TextCell $q3.Range("C5") "创金合信量化多因子股票C"
Set-TextCell $q3.Range("D5") "0.75"
Set-TextCell $q3.Range("E5") "91.71"
Set-TextCell $q3.Range("F5") "1.22"
Set-TextCell $q3.Range("G5") "0.0092"
$q3.Range("H5").Value = 10

# ---------------------------------------------------------------------
# 2) Insert a new row 2 in "总计" for the 2022-Q3 totals, pushing the
#    existing rows down by one (their date labels shift by one quarter).
#    The "A" column is a plain 0-based row index, so every shifted row
#    needs its index bumped by one as well.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Re-apply the bold/bordered index-column style (copied from the row
# that just shifted down to row 3) to the newly inserted A2 cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.07000000000000001

# Renumber the index column for the rows that shifted down (old index
# N now sits one row lower and must read N+1).
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7
